$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 21:12:32"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 21:12:28"
$wsZhCn.Range("K2").Value = "2016-09-04 21:12:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 21:12:32"
$wsDeDe.Range("K2").Value = "2016-09-04 21:12:54"
